$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the price of the Wooden Marble Labyrinth (row 2)
$ws.Range("C2").Value = 20.16

# Re-enter the Extension formula across the whole range so Excel
# collapses it into a shared formula group (E2:E9)
$ws.Range("E2:E9").Formula = "=D2*C2"

# Update the selected cell/range
$ws.Range("C3").Select()
